$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.065.03"
$ws.Range("E2").Value = "  +3.50%  "
$ws.Range("D3").Value = "3.820.98"
$ws.Range("E3").Value = "  +7.86%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "428.93"
$ws.Range("E5").Value = "  +8.47%  "
$ws.Range("D6").Value = "131.45"
$ws.Range("E6").Value = "  +3.66%  "
$ws.Range("D7").Value = "3.811.88"
$ws.Range("E7").Value = "  +7.87%  "
$ws.Range("E8").Value = "  +3.48%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "0.736"
$ws.Range("E10").Value = "  +6.70%  "
$ws.Range("E11").Value = "  +4.14%  "
$ws.Range("D12").Value = "0.0000338"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").Value = "41.64"
$ws.Range("E13").Value = "  +5.80%  "
$ws.Range("D14").Value = "10.53"
$ws.Range("E14").Value = "  +13.24%  "
$ws.Range("D15").Value = "4.428.16"
$ws.Range("E15").Value = "  +8.08%  "
$ws.Range("D16").Value = "15.41"
$ws.Range("E16").Value = "  +20.92%  "
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("D18").Value = "3.809.22"
$ws.Range("E18").Value = "  +7.77%  "
$ws.Range("E19").Value = "  +6.40%  "
$ws.Range("D20").Value = "1.12"
$ws.Range("E20").Value = "  +8.23%  "
$ws.Range("D21").Value = "66.264.94"
$ws.Range("E21").Value = "  +3.72%  "
$ws.Range("D22").Value = "416.13"
$ws.Range("E22").Value = "  +3.97%  "
$ws.Range("D23").Value = "15.23"
$ws.Range("E23").Value = "  +8.73%  "
$ws.Range("D24").Value = "85.57"
$ws.Range("E24").Value = "  +4.50%  "
$ws.Range("E25").Value = "  +8.24%  "
$ws.Range("D26").Value = "37.21"
$ws.Range("E26").Value = "  +8.45%  "
$ws.Range("D27").Value = "10.10"
$ws.Range("E27").Value = "  +13.72%  "
$ws.Range("D28").Value = "3.30"
$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").Value = "5.38"
$ws.Range("E29").Value = "  +2.61%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "9.34"
$ws.Range("E30").Value = "  +35.78%  "
$ws.Range("D31").Value = "14.19"
$ws.Range("E31").Value = "  +18.75%  "
$ws.Range("D32").Value = "710.50"
$ws.Range("E32").Value = "  +4.62%  "
$ws.Range("E33").Value = "  +13.16%  "
$ws.Range("E34").Value = "  +7.38%  "
$ws.Range("D35").Value = "5.83"
$ws.Range("E35").Value = "  +40.95%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "38.98"
$ws.Range("E37").Value = "  +4.77%  "
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("D39").Value = "55.79"
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("D40").Value = "0.0471"
$ws.Range("E40").Value = "  +6.81%  "
$ws.Range("D41").Value = "0.0₃0729"
$ws.Range("E41").Value = "  +16.99%  "
$ws.Range("D42").Value = "2.89"
$ws.Range("E42").Value = "  +2.60%  "
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("E44").Value = "  +4.45%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "3.25"
$ws.Range("E45").Value = "  +5.98%  "
$ws.Range("B46").Value = "LidoDAOToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D46").Value = "3.40"
$ws.Range("E46").Value = "  +9.48%  "
$ws.Range("E47").Value = "  +16.89%  "
$ws.Range("E48").Value = "  +43.02%  "
$ws.Range("E49").Value = "  +5.66%  "
$ws.Range("E50").Value = "  +5.22%  "
$ws.Range("D51").Value = "2.84"
$ws.Range("E51").Value = "  +4.17%  "
